$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(5).Insert()
$ws.Rows.Item(10).Insert()

$ws.Rows.Item(4).Style = "Normal"
